$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.044.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.277.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "272.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.96%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.620.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.833"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.262.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.036.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0908"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +20.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.249"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.432"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.504.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.41%  "
